{"js": "// Apply the \"Trust Before Intelligence One Pager\" bio update:\n//   - Expand \"Ram Katamaraja\" to \"Ram Dhan Yadav Katamaraja\" (byline + bio paragraph)\n//   - Reorder / reword the awards sentence in the bio paragraph\n//   - Update \"and MIT.\" -> \"MIT to name a few.\" and \"Ram is\" -> \"Ram Dhan is\"\n//   - \"a decade of enterprise AI implementations across\" -> \"three decades of Enterprise architecture work in\"\n//   - \"agriculture and utility services.\" -> \"agriculture, utilities, banking, finance, supply chain, telecom industries.\"\n\nasync function replaceOnce(context, findText, replaceText) {\n  const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Title byline: \"By Ram Katamaraja\" -> \"By Ram Dhan Yadav Katamaraja\"\nawait replaceOnce(context, \"By Ram Katamaraja\", \"By Ram Dhan Yadav Katamaraja\");\n\n// 2. Bio paragraph, opening name\nawait replaceOnce(context, \"Ram Katamaraja founded\", \"Ram Dhan Yadav Katamaraja founded\");\n\n// 3. Bio paragraph, swap the two award clauses and drop the word \"Prize\"\nawait replaceOnce(\n  context,\n  \"He received the McGovern Foundation's \\\"AI for the Betterment of Humanity Prize\\\" and MIT's Most Promising Work of the Future Solution awards.\",\n  \"He received the MIT's Most Promising Work of the Future Solution and McGovern Foundation's \\\"AI for the Betterment of Humanity\\\" awards.\"\n);\n\n// 4. Bio paragraph, \"and MIT.\" -> \"MIT to name a few.\" and \"Ram is\" -> \"Ram Dhan is\"\nawait replaceOnce(\n  context,\n  \"Harvard Business School, and MIT. Ram is an alumn\",\n  \"Harvard Business School, MIT to name a few. Ram Dhan is an alumn\"\n);\n\n// 5. Bio paragraph, \"a decade of enterprise AI implementations across\" -> \"three decades of Enterprise architecture work in\"\nawait replaceOnce(\n  context,\n  \"His frameworks emerged from a decade of enterprise AI implementations across\",\n  \"His frameworks emerged from three decades of Enterprise architecture work in\"\n);\n\n// 6. Bio paragraph, expand the industries list\nawait replaceOnce(\n  context,\n  \"agriculture and utility services.\",\n  \"agriculture, utilities, banking, finance, supply chain, telecom industries.\"\n);\n", "ps1": "# Apply the \"Trust Before Intelligence One Pager\" bio update:\n#   - Expand \"Ram Katamaraja\" to \"Ram Dhan Yadav Katamaraja\" (byline + bio paragraph)\n#   - Reorder / reword the awards sentence in the bio paragraph\n#   - Update \"and MIT.\" -> \"MIT to name a few.\" and \"Ram is\" -> \"Ram Dhan is\"\n#   - \"a decade of enterprise AI implementations across\" -> \"three decades of Enterprise architecture work in\"\n#   - \"agriculture and utility services.\" -> \"agriculture, utilities, banking, finance, supply chain, telecom industries.\"\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $findText\n    $range.Find.Replacement.Text = $replaceText\n    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\n# 1. Title byline\nReplace-Text \"By Ram Katamaraja\" \"By Ram Dhan Yadav Katamaraja\"\n\n# 2. Bio paragraph, opening name\nReplace-Text \"Ram Katamaraja founded\" \"Ram Dhan Yadav Katamaraja founded\"\n\n# 3. Bio paragraph, swap the two award clauses and drop the word \"Prize\"\nReplace-Text \"He received the McGovern Foundation's `\"AI for the Betterment of Humanity Prize`\" and MIT's Most Promising Work of the Future Solution awards.\" \"He received the MIT's Most Promising Work of the Future Solution and McGovern Foundation's `\"AI for the Betterment of Humanity`\" awards.\"\n\n# 4. Bio paragraph, \"and MIT.\" -> \"MIT to name a few.\" and \"Ram is\" -> \"Ram Dhan is\"\nReplace-Text \"Harvard Business School, and MIT. Ram is an alumn\" \"Harvard Business School, MIT to name a few. Ram Dhan is an alumn\"\n\n# 5. Bio paragraph, \"a decade of enterprise AI implementations across\" -> \"three decades of Enterprise architecture work in\"\nReplace-Text \"His frameworks emerged from a decade of enterprise AI implementations across\" \"His frameworks emerged from three decades of Enterprise architecture work in\"\n\n# 6. Bio paragraph, expand the industries list\nReplace-Text \"agriculture and utility services.\" \"agriculture, utilities, banking, finance, supply chain, telecom industries.\"\n"}
